$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be parsed as a number
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated price / volume values
$ws.Range("D2").Value = "22.117.93"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.556.23"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "0.9990"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "288.00"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "0.3793"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("D8").Value = "0.3290"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "43.28"
$ws.Range("E9").Value = "  -10.50%  "
$ws.Range("D10").Value = "1.140"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "0.07369"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "20.18"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "6.834"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "1.551.27"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "0.00001105"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "0.06620"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "85.94"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "6.398"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "11.71"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "22.110.49"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").Value = "2.528"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").Value = "150.57"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "19.14"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "4.913"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").Value = "121.78"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("D31").Value = "1.726.70"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "1.081"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "5.965"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").Value = "1.851"
$ws.Range("E34").Value = "  -8.04%  "
$ws.Range("D35").Value = "0.08231"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "9.338"
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("D37").Value = "0.02336"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D38").Value = "0.06259"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "5.298"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "0.2164"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "1.256"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "11.06"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").Value = "0.6064"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "13.81"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "3.738"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "0.5860"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("D48").Value = "1.992"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").Value = "122.40"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "1.178"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "0.07023"
$ws.Range("E51").Value = "  -2.78%  "
